$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in uncertainty values (column C) for rows 27-30
$ws.Range("C27").Value = 0.05
$ws.Range("C28").Value = 0.04
$ws.Range("C29").Value = 0.05
$ws.Range("C30").Value = 0.04

# Update the view: scroll so row 25 is at the top, and move the
# active selection to E28 (matches the saved sheetView/selection state)
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E28").Select()
